$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, pushing the existing row 22 (and below)
# down by one. This preserves the old row 22 -> new row 23 and old row 23
# -> new row 24. Inserting via the full row-width range (rather than the
# whole Rows collection) carries the existing formatting down automatically.
$ws.Range("A22:R22").Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Range("A22").Value = 11
$ws.Range("B22").Value = "Vega Monumental Concepción"
$ws.Range("C22").Value = "Bíobío"
$ws.Range("D22").Value = 44615
$ws.Range("E22").Value = 8
$ws.Range("F22").Value = 100112022
$ws.Range("G22").Value = "Arveja Verde"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 28000
$ws.Range("L22").Value = 30000
$ws.Range("M22").Value = 29000
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Carahue"
$ws.Range("P22").Value = 1160
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
